$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 189.29411
$ws.Range("I39").Value = 16.583334
$ws.Range("J39").Value = 603.8
$ws.Range("K39").Value = 49.750002
$ws.Range("L39").Value = 1811.4
$ws.Range("M39").Value = 246.249998
$ws.Range("N39").Value = -2403.4
$ws.Range("H51").Value = 2807.5386
$ws.Range("J51").Value = 2999.75
$ws.Range("L51").Value = 2999.75
$ws.Range("N51").Value = -3967.75
$ws.Range("H96").Value = 844181.7
$ws.Range("J96").Value = 25291.8
$ws.Range("L96").Value = 75875.39999999999
$ws.Range("N96").Value = -78621.39999999999
$ws.Range("H131").Value = 7351.769
$ws.Range("I131").Value = 6321.875
$ws.Range("K131").Value = 18965.625
$ws.Range("M131").Value = -13925.625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 8000.8335
$ws.Range("I63").Value = 4002.5
$ws.Range("K63").Value = 4002.5
$ws.Range("M63").Value = -3316.5
$ws.Range("H66").Value = 8000.8335
$ws.Range("I66").Value = 4002.5
$ws.Range("K66").Value = 20012.5
$ws.Range("M66").Value = -16580.5
$ws.Range("H108").Value = 20000
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H132").Value = 2071.5417
$ws.Range("I132").Value = 2071.5417
$ws.Range("K132").Value = 6214.625100000001
$ws.Range("M132").Value = -3684.625100000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4696.923
$ws.Range("I31").Value = 2418.0386
$ws.Range("K31").Value = 2418.0386
$ws.Range("M31").Value = -2123.0386
$ws.Range("H34").Value = 4696.923
$ws.Range("I34").Value = 2418.0386
$ws.Range("K34").Value = 2418.0386
$ws.Range("M34").Value = -2216.0386
$ws.Range("H93").Value = 33579.6
$ws.Range("I93").Value = 29474.5
$ws.Range("K93").Value = 29474.5
$ws.Range("M93").Value = -27602.5
$ws.Range("H105").Value = 1246.25
$ws.Range("I105").Value = 1331.6666
$ws.Range("J105").Value = 990
$ws.Range("K105").Value = 1331.6666
$ws.Range("L105").Value = 990
$ws.Range("M105").Value = 415.3334
$ws.Range("N105").Value = -4484
$ws.Range("H132").Value = 2682.1516
$ws.Range("I132").Value = 2770.4
$ws.Range("K132").Value = 8311.200000000001
$ws.Range("M132").Value = -5781.200000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 38
$ws.Range("I15").Value = 40
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 120
$ws.Range("L15").Value = 90
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -370
$ws.Range("H26").Value = 624.8889
$ws.Range("J26").Value = 200
$ws.Range("L26").Value = 600
$ws.Range("N26").Value = -1176
$ws.Range("H32").Value = 2750
$ws.Range("J32").Value = 2750
$ws.Range("L32").Value = 8250
$ws.Range("N32").Value = -8816
$ws.Range("H34").Value = 625
$ws.Range("I34").Value = 735.7143
$ws.Range("J34").Value = 366.66666
$ws.Range("K34").Value = 2207.1429
$ws.Range("L34").Value = 1099.99998
$ws.Range("M34").Value = -2123.1429
$ws.Range("N34").Value = -1267.99998
$ws.Range("H46").Value = 9967436
$ws.Range("I46").Value = 23232972
$ws.Range("J46").Value = 3334668.8
$ws.Range("K46").Value = 69698916
$ws.Range("L46").Value = 10004006.4
$ws.Range("M46").Value = -69698825
$ws.Range("N46").Value = -10004188.4
$ws.Range("H49").Value = 1001.5
$ws.Range("I49").Value = 1003
$ws.Range("J49").Value = 1000
$ws.Range("K49").Value = 3009
$ws.Range("L49").Value = 3000
$ws.Range("M49").Value = -2853
$ws.Range("N49").Value = -3312
$ws.Range("H63").Value = 950
$ws.Range("I63").Value = 950
$ws.Range("K63").Value = 2850
$ws.Range("M63").Value = -2101
$ws.Range("H66").Value = 950
$ws.Range("I66").Value = 950
$ws.Range("K66").Value = 8550
$ws.Range("M66").Value = -4806
$ws.Range("H74").Value = 20541
$ws.Range("J74").Value = 25600
$ws.Range("L74").Value = 76800
$ws.Range("N74").Value = -78922
$ws.Range("H77").Value = 20541
$ws.Range("J77").Value = 25600
$ws.Range("L77").Value = 230400
$ws.Range("N77").Value = -241008
$ws.Range("H97").Value = 314.3
$ws.Range("J97").Value = 340.6875
$ws.Range("L97").Value = 1022.0625
$ws.Range("N97").Value = -2014.0625
$ws.Range("H137").Value = 4797.1562
$ws.Range("I137").Value = 2064.625
$ws.Range("J137").Value = 5708
$ws.Range("K137").Value = 6193.875
$ws.Range("L137").Value = 17124
$ws.Range("M137").Value = -1093.875
$ws.Range("N137").Value = -27324
$ws.Range("H141").Value = 3999
$ws.Range("I141").Value = 3999
$ws.Range("K141").Value = 11997
$ws.Range("M141").Value = -6817

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 155.4
$ws.Range("I2").Value = 58.142857
$ws.Range("J2").Value = 240.5
$ws.Range("K2").Value = 58.142857
$ws.Range("L2").Value = 240.5
$ws.Range("M2").Value = 54.857143
$ws.Range("N2").Value = -466.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5379.8
$ws.Range("I7").Value = 4849.75
$ws.Range("K7").Value = 4849.75
$ws.Range("M7").Value = -4737.75
$ws.Range("H22").Value = 3428.889
$ws.Range("I22").Value = 3461.4285
$ws.Range("J22").Value = 3408.182
$ws.Range("K22").Value = 3461.4285
$ws.Range("L22").Value = 3408.182
$ws.Range("M22").Value = -3166.4285
$ws.Range("N22").Value = -3998.182
$ws.Range("H27").Value = 3428.889
$ws.Range("I27").Value = 3461.4285
$ws.Range("J27").Value = 3408.182
$ws.Range("K27").Value = 3461.4285
$ws.Range("L27").Value = 3408.182
$ws.Range("M27").Value = -3354.4285
$ws.Range("N27").Value = -3622.182
$ws.Range("H126").Value = 5379.8
$ws.Range("I126").Value = 4849.75
$ws.Range("K126").Value = 14549.25
$ws.Range("M126").Value = -12079.25
$ws.Range("H132").Value = 4218.1816
$ws.Range("I132").Value = 5080
$ws.Range("K132").Value = 15240
$ws.Range("M132").Value = -12710
$ws.Range("H136").Value = 3473.8774
$ws.Range("J136").Value = 3227.7576
$ws.Range("L136").Value = 9683.272799999999
$ws.Range("N136").Value = -14783.2728

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 609.1177
$ws.Range("I100").Value = 564.93335
$ws.Range("K100").Value = 1129.8667
$ws.Range("M100").Value = -588.8667
$ws.Range("H122").Value = 6155.778
$ws.Range("I122").Value = 2257.7144
$ws.Range("K122").Value = 6773.1432
$ws.Range("M122").Value = -4323.1432
$ws.Range("H126").Value = 1311.909
$ws.Range("J126").Value = 1902.5
$ws.Range("L126").Value = 5707.5
$ws.Range("N126").Value = -10647.5
$ws.Range("H132").Value = 1888.6471
$ws.Range("I132").Value = 1682.3636
$ws.Range("J132").Value = 2266.8333
$ws.Range("K132").Value = 5047.0908
$ws.Range("L132").Value = 6800.499899999999
$ws.Range("M132").Value = -2517.0908
$ws.Range("N132").Value = -11860.4999
$ws.Range("H136").Value = 1063.4
$ws.Range("I136").Value = 860.0909
$ws.Range("J136").Value = 1622.5
$ws.Range("K136").Value = 2580.2727
$ws.Range("L136").Value = 4867.5
$ws.Range("M136").Value = -30.27269999999999
$ws.Range("N136").Value = -9967.5

Write-Host "Applied all Goblin_Profits updates"